$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 652 (shifts existing rows 652:749 down to 653:750,
# new row inherits formatting - e.g. the date-formatted style - from the row above).
$ws.Rows("652:652").Insert()

# Populate the newly inserted row with the new price record.
$ws.Range("A652").Value = 11
$ws.Range("B652").Value = "Vega Monumental Concepción"
$ws.Range("C652").Value = "Bíobío"
$ws.Range("D652").Value = "2023-04-18"
$ws.Range("E652").Value = 8
$ws.Range("F652").Value = "Fruta"
$ws.Range("G652").Value = 100108
$ws.Range("H652").Value = "Tropicales y subtropicales"
$ws.Range("I652").Value = 100108006
$ws.Range("J652").Value = "Plátano"
$ws.Range("K652").Value = "Sin especificar"
$ws.Range("L652").Value = "Pintón"
$ws.Range("M652").Value = 1100
$ws.Range("N652").Value = 18000
$ws.Range("O652").Value = 20000
$ws.Range("P652").Value = 18909
$ws.Range("Q652").Value = "$/caja 20 kilos"
$ws.Range("R652").Value = "Ecuador"
$ws.Range("S652").Value = 945
$ws.Range("T652").Value = 20
